# Convert the two Word field constructs ({ m: userdoc 'zone1' } and
# { m:enduserdoc }) that are currently stored as real field codes
# (w:fldChar/w:instrText) into plain literal text runs (w:t), while
# keeping the _GoBack bookmark that sits inside the first field.
#
# Strategy: for each field paragraph we
#   1) insert the replacement run(s) (as raw WordprocessingML) right
#      before the field, via a collapsed Range + InsertXML - this is a
#      pure insertion, it does not disturb the existing field runs;
#   2) then remove the original field (begin/instrText*/end) with
#      Field.Delete(), which cleanly deletes exactly the field's runs.

$d = $word.ActiveDocument

# --- Paragraph 2: "{m: userdoc 'zone1' }" field -> literal runs,
#     keeping the _GoBack bookmark in place between the space-run and
#     the "userdoc 'zone1'" run. ---
$p2 = $d.Paragraphs.Item(2)
$insertPoint2 = $d.Range($p2.Range.Start, $p2.Range.Start)
$xml2 = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:r><w:t>{</w:t></w:r>
<w:r><w:t>m</w:t></w:r>
<w:r><w:t>:</w:t></w:r>
<w:r><w:t xml:space='preserve'> </w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
<w:r><w:t>userdoc 'zone1'</w:t></w:r>
<w:r><w:t xml:space='preserve'>}</w:t></w:r>
</w:p>
"@
$insertPoint2.InsertXML($xml2)
$d.Fields.Item(1).Delete()

# --- Paragraph 4: "{m:enduserdoc}" field -> a single literal run. ---
$p4 = $d.Paragraphs.Item(4)
$insertPoint4 = $d.Range($p4.Range.Start, $p4.Range.Start)
$xml4 = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:r><w:t xml:space='preserve'>{m:enduserdoc}</w:t></w:r>
</w:p>
"@
$insertPoint4.InsertXML($xml4)
$d.Fields.Item(1).Delete()
